$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2409.4
$ws.Range("I32").Value = 1465
$ws.Range("J32").Value = 2814.1428
$ws.Range("K32").Value = 1465
$ws.Range("L32").Value = 2814.1428
$ws.Range("M32").Value = -1139
$ws.Range("N32").Value = -3466.1428
$ws.Range("H40").Value = 4000
$ws.Range("J40").Value = 4000
$ws.Range("L40").Value = 4000
$ws.Range("N40").Value = -4350
$ws.Range("H64").Value = 6749
$ws.Range("H67").Value = 6749
$ws.Range("H69").Value = 9824.081
$ws.Range("I69").Value = 8605.4
$ws.Range("K69").Value = 25816.2
$ws.Range("M69").Value = -24942.2
$ws.Range("H72").Value = 9824.081
$ws.Range("I72").Value = 8605.4
$ws.Range("K72").Value = 77448.59999999999
$ws.Range("M72").Value = -73080.59999999999
$ws.Range("H97").Value = 1881.5
$ws.Range("J97").Value = 1881.5
$ws.Range("L97").Value = 5644.5
$ws.Range("N97").Value = -6636.5
$ws.Range("H112").Value = 2592.9312
$ws.Range("I112").Value = 2327.0908
$ws.Range("J112").Value = 2755.389
$ws.Range("K112").Value = 6981.2724
$ws.Range("L112").Value = 8266.167000000001
$ws.Range("M112").Value = -5873.2724
$ws.Range("N112").Value = -10482.167
$ws.Range("H133").Value = 84800
$ws.Range("J133").Value = 84800
$ws.Range("L133").Value = 84800
$ws.Range("N133").Value = -94920
$ws.Range("H138").Value = 19609844
$ws.Range("I138").Value = 1252.8148
$ws.Range("J138").Value = 41669508
$ws.Range("K138").Value = 3758.4444
$ws.Range("L138").Value = 125008524
$ws.Range("M138").Value = 1381.5556
$ws.Range("N138").Value = -125018804

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 53550
$ws.Range("J139").Value = 76000
$ws.Range("L139").Value = 76000
$ws.Range("N139").Value = -86280

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 979.86664
$ws.Range("I94").Value = 964.8
$ws.Range("K94").Value = 964.8
$ws.Range("M94").Value = -513.8
$ws.Range("H96").Value = 55397.4
$ws.Range("J96").Value = 93994
$ws.Range("L96").Value = 93994
$ws.Range("N96").Value = -99486
$ws.Range("H99").Value = 2285.913
$ws.Range("I99").Value = 1995.8889
$ws.Range("K99").Value = 1995.8889
$ws.Range("M99").Value = -497.8888999999999
$ws.Range("H107").Value = 3479.3635
$ws.Range("J107").Value = 3814.3333
$ws.Range("L107").Value = 3814.3333
$ws.Range("N107").Value = -7654.3333
$ws.Range("H134").Value = 3512101.8
$ws.Range("I134").Value = 3924613.8
$ws.Range("J134").Value = 5750
$ws.Range("K134").Value = 11773841.4
$ws.Range("L134").Value = 17250
$ws.Range("M134").Value = -11771306.4
$ws.Range("N134").Value = -22320

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9826.333000000001
$ws.Range("J31").Value = 10751.387
$ws.Range("L31").Value = 10751.387
$ws.Range("N31").Value = -11341.387
$ws.Range("H34").Value = 9826.333000000001
$ws.Range("J34").Value = 10751.387
$ws.Range("L34").Value = 10751.387
$ws.Range("N34").Value = -11155.387
$ws.Range("H58").Value = 2699.4355
$ws.Range("I58").Value = 2593.827
$ws.Range("J58").Value = 3248.6
$ws.Range("K58").Value = 2593.827
$ws.Range("L58").Value = 3248.6
$ws.Range("M58").Value = -2390.827
$ws.Range("N58").Value = -3654.6
$ws.Range("H87").Value = 85575.39999999999
$ws.Range("J87").Value = 85575.39999999999
$ws.Range("L87").Value = 85575.39999999999
$ws.Range("N87").Value = -87947.39999999999
$ws.Range("H90").Value = 85575.39999999999
$ws.Range("J90").Value = 85575.39999999999
$ws.Range("L90").Value = 256726.2
$ws.Range("N90").Value = -268582.2
$ws.Range("H114").Value = 69994.25
$ws.Range("J114").Value = 69994.25
$ws.Range("L114").Value = 69994.25
$ws.Range("N114").Value = -78672.25
$ws.Range("H129").Value = 59728
$ws.Range("J129").Value = 59728
$ws.Range("L129").Value = 59728
$ws.Range("N129").Value = -69728
$ws.Range("H136").Value = 2699.4355
$ws.Range("I136").Value = 2593.827
$ws.Range("J136").Value = 3248.6
$ws.Range("K136").Value = 7781.481000000001
$ws.Range("L136").Value = 9745.799999999999
$ws.Range("M136").Value = -5231.481000000001
$ws.Range("N136").Value = -14845.8

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2373.1875
$ws.Range("I109").Value = 1297.1
$ws.Range("J109").Value = 4166.6665
$ws.Range("K109").Value = 3891.3
$ws.Range("L109").Value = 12499.9995
$ws.Range("M109").Value = -2851.3
$ws.Range("N109").Value = -14579.9995
$ws.Range("H131").Value = 1604.9706
$ws.Range("J131").Value = 1705.0385
$ws.Range("L131").Value = 5115.1155
$ws.Range("N131").Value = -15195.1155
$ws.Range("H132").Value = 1001315.8
$ws.Range("J132").Value = 1001315.8
$ws.Range("L132").Value = 9011842.200000001
$ws.Range("N132").Value = -9016902.200000001
$ws.Range("H134").Value = 2751.1667
$ws.Range("I134").Value = 2751.1667
$ws.Range("K134").Value = 8253.500100000001
$ws.Range("M134").Value = -3183.500100000001
$ws.Range("H140").Value = 1944.5454
$ws.Range("I140").Value = 1462.1052
$ws.Range("K140").Value = 4386.3156
$ws.Range("M140").Value = 793.6844000000001

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 949.8182
$ws.Range("I97").Value = 666.44446
$ws.Range("J97").Value = 2225
$ws.Range("K97").Value = 666.44446
$ws.Range("L97").Value = 2225
$ws.Range("M97").Value = -170.44446
$ws.Range("N97").Value = -3217

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 14448.4375
$ws.Range("I122").Value = 15746.75
$ws.Range("J122").Value = 12284.583
$ws.Range("K122").Value = 47240.25
$ws.Range("L122").Value = 36853.749
$ws.Range("M122").Value = -44790.25
$ws.Range("N122").Value = -41753.749
$ws.Range("H132").Value = 5462.5
$ws.Range("I132").Value = 5316.6665
$ws.Range("J132").Value = 5900
$ws.Range("K132").Value = 15949.9995
$ws.Range("L132").Value = 17700
$ws.Range("M132").Value = -13419.9995
$ws.Range("N132").Value = -22760
$ws.Range("H133").Value = 29888
$ws.Range("J133").Value = 29888
$ws.Range("L133").Value = 29888
$ws.Range("N133").Value = -34948
$ws.Range("H136").Value = 6085
$ws.Range("J136").Value = 8449.6
$ws.Range("L136").Value = 25348.8
$ws.Range("N136").Value = -30448.8

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 47513.5
$ws.Range("J51").Value = 59998
$ws.Range("L51").Value = 59998
$ws.Range("N51").Value = -61018
$ws.Range("H123").Value = 57865.8
$ws.Range("J123").Value = 63109.668
$ws.Range("L123").Value = 63109.668
$ws.Range("N123").Value = -72909.66800000001
$ws.Range("H126").Value = 1902.3158
$ws.Range("I126").Value = 1230.3077
$ws.Range("J126").Value = 3358.3333
$ws.Range("K126").Value = 3690.9231
$ws.Range("L126").Value = 10074.9999
$ws.Range("M126").Value = -1220.9231
$ws.Range("N126").Value = -15014.9999
$ws.Range("H132").Value = 3907.6155
$ws.Range("I132").Value = 3887.5
$ws.Range("J132").Value = 3939.8
$ws.Range("K132").Value = 11662.5
$ws.Range("L132").Value = 11819.4
$ws.Range("M132").Value = -9132.5
$ws.Range("N132").Value = -16879.4
$ws.Range("H136").Value = 1878.5
$ws.Range("I136").Value = 1745
$ws.Range("J136").Value = 2088.2856
$ws.Range("K136").Value = 5235
$ws.Range("L136").Value = 6264.8568
$ws.Range("M136").Value = -2685
$ws.Range("N136").Value = -11364.8568
